# Update database: roll the 5-year reporting window forward by one year.
# Drop the oldest period (1396/12) and add the newest period (1401/12),
# shifting every existing column one slot to the left (D<-E, E<-F, F<-G, G<-H)
# and filling column H with the newly published figures.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-Row($row, $values) {
    $cols = @("D", "E", "F", "G", "H")
    for ($i = 0; $i -lt 5; $i++) {
        $ws.Range($cols[$i] + $row).Value = $values[$i]
    }
}

# Row 8: financial-period headers (دوره مالی)
Set-Row 8 @(
    "12 ماهه منتهی به 1397/12",
    "12 ماهه منتهی به 1398/12",
    "12 ماهه منتهی به 1399/12",
    "12 ماهه منتهی به 1400/12",
    "12 ماهه منتهی به 1401/12"
)

# Row 9: publish-date headers (تاریخ انتشار)
# H9 is a bare "YYYY-MM-DD" string; force text formatting so it is not
# auto-converted into a date serial number.
$ws.Range("H9").NumberFormat = "@"
Set-Row 9 @(
    "1399-03-13 (8)",
    "1400-03-11 (10)",
    "1401-03-24 (10)",
    "1402-02-28 (8)",
    "1402-02-28"
)

# Row 11: فروش (Sales)
Set-Row 11 @(20664, 27711, 32244, 36259, 37078)

# Row 12: بهای تمام شده کالای فروش رفته (Cost of goods sold)
Set-Row 12 @(-11915, -15394, -14319, -16711, -19748)

# Row 13: سود (زیان) ناخالص (Gross profit)
Set-Row 13 @(8748, 12317, 17925, 19547, 17330)

# Row 14: هزینه های عمومی, اداری و تشکیلاتی (G&A expenses)
Set-Row 14 @(-1611, -1781, -1478, -1729, -1690)

# Row 16: خالص سایر درامدها (هزینه ها) ی عملیاتی (Other operating income/expense, net)
Set-Row 16 @(421, 1054, 597, 122, 932)

# Row 17: سود (زیان) عملیاتی (Operating profit)
Set-Row 17 @(7559, 11590, 17043, 17940, 16572)

# Row 18: هزینه های مالی (Financial expenses)
Set-Row 18 @(-1728, -1401, -313, -555, -1308)

# Row 19: خالص سایر درامدها و هزینه های غیرعملیاتی (Other non-operating income/expense, net)
Set-Row 19 @(79, -519, 229, 307, 137)

# Row 20: سود (زیان) خالص عملیات در حال تداوم قبل از مالیات (Pre-tax profit from continuing operations)
Set-Row 20 @(5910, 9671, 16960, 17693, 15401)

# Row 21: مالیات (Tax)
Set-Row 21 @(-465, -939, -1068, -2183, -2145)

# Row 22: سود (زیان) خالص عملیات در حال تداوم (Net profit from continuing operations)
Set-Row 22 @(5445, 8731, 15892, 15509, 13256)

# Row 24: سود (زیان) خالص (Net profit)
Set-Row 24 @(5445, 8731, 15892, 15509, 13256)

# Row 26: سرمایه (Capital)
Set-Row 26 @(4943, 3898, 3162, 2710, 2026)
